# edit.ps1 - applies the "falcon-7b-instruct / summary & visualisation" rework
# to the single-slide Abu Dhabi Open Data Intelligence workflow deck.

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

function Get-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Name -eq $name) {
            return $sh
        }
    }
    return $null
}

# Shape.Left/Top/Width/Height round-trip through a 32-bit "points" float
# (exactly like real PowerPoint's Single-typed properties), so naive
# EMU/12700 assignments can truncate to the EMU value below the true
# target. Search nearby float32 values until the stored EMU (computed the
# same way the host does: truncate(float32(points) * 12700)) matches
# exactly, so the saved XML reproduces the target EMU bit-for-bit.
function EmuToPtPrecise($targetEmu) {
    $base = $targetEmu / 12700.0
    for ($i = 0; $i -lt 2000; $i++) {
        $cand = $base + ($i * 0.000001)
        $f32 = [float]$cand
        $emu = [int]([double]$f32 * 12700.0)
        if ($emu -eq $targetEmu) {
            return $cand
        }
    }
    return $base
}

function Set-ShapeGeometryEmu($shape, $left, $top, $width, $height) {
    if ($null -ne $left)   { $shape.Left   = EmuToPtPrecise $left }
    if ($null -ne $top)    { $shape.Top    = EmuToPtPrecise $top }
    if ($null -ne $width)  { $shape.Width  = EmuToPtPrecise $width }
    if ($null -ne $height) { $shape.Height = EmuToPtPrecise $height }
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# 1) "Web " + "Application (" -> single run "Web Application ("
# ---------------------------------------------------------------------------

$webAppBox = Get-ShapeByName $s "TextBox 24"
$webAppRange = $webAppBox.TextFrame.TextRange.Paragraphs(1, 1)
$webAppRange.Characters(1, 17).Text = "Web Application ("

# ---------------------------------------------------------------------------
# 2) "Using LLM (... tiiuae/falcon-7b)" -> "...tiiuae/falcon-7b-instruct)"
#    box also grows taller (738664 -> 954107 EMU)
# ---------------------------------------------------------------------------

$usingLlmBox = Get-ShapeByName $s "TextBox 54"
$usingLlmRange = $usingLlmBox.TextFrame.TextRange
$start = $usingLlmRange.Text.IndexOf("/falcon-7b") + 1
$usingLlmRange.Characters($start, 10).Text = "/falcon-7b-instruct"
Set-ShapeGeometryEmu $usingLlmBox $null $null $null 954107

# ---------------------------------------------------------------------------
# 3) "Train data categories using tiiuae/Falcon-7b" -> "...Falcon-7b-instruct"
#    box moves up and grows taller
# ---------------------------------------------------------------------------

$trainDataBox = Get-ShapeByName $s "TextBox 72"
$trainDataPara2 = $trainDataBox.TextFrame.TextRange.Paragraphs(2, 1)
$start = $trainDataPara2.Text.IndexOf("/Falcon-7b") + 1
$trainDataPara2.Characters($start, 10).Text = "/Falcon-7b-instruct"
Set-ShapeGeometryEmu $trainDataBox $null 4647681 $null 1169551

# ---------------------------------------------------------------------------
# 4) Connector feeding "Train Data" box: drop the vertical flip, stretch it
#    to follow the box's new position
# ---------------------------------------------------------------------------

$trainDataConn = Get-ShapeByName $s "Straight Arrow Connector 76"
Set-ShapeGeometryEmu $trainDataConn $null 5229240 $null 3217
$trainDataConn.VerticalFlip = 0

# ---------------------------------------------------------------------------
# 5) "Response in text" -> "Summary & " + "Visualisation"
#    box moves up and grows taller
# ---------------------------------------------------------------------------

$outputBox = Get-ShapeByName $s "TextBox 92"
$outputPara2 = $outputBox.TextFrame.TextRange.Paragraphs(2, 1)
$outputPara2.Characters(1, $outputPara2.Length).Text = "Summary & Visualisation"
$outputPara2b = $outputBox.TextFrame.TextRange.Paragraphs(2, 1)
$outputPara2b.Characters(1, 10).Text = "Summary & "
Set-ShapeGeometryEmu $outputBox $null 3397018 $null 738664

# ---------------------------------------------------------------------------
# 6) Connector feeding the "Convert ... to audio" box: stretch it downward
# ---------------------------------------------------------------------------

$outputConn = Get-ShapeByName $s "Straight Arrow Connector 94"
Set-ShapeGeometryEmu $outputConn $null $null $null 3218

# ---------------------------------------------------------------------------
# 7) "Convert text to audio" -> "Convert " + "summary " + "to audio"
#    box grows taller
# ---------------------------------------------------------------------------

$audioBox = Get-ShapeByName $s "TextBox 96"
$audioPara2 = $audioBox.TextFrame.TextRange.Paragraphs(2, 1)
$audioPara2.Characters(1, 8).Text = "Convert "
$audioPara2b = $audioBox.TextFrame.TextRange.Paragraphs(2, 1)
$audioPara2b.Characters(9, 5).Text = "summary "
$audioPara2c = $audioBox.TextFrame.TextRange.Paragraphs(2, 1)
$audioPara2c.Characters(17, $audioPara2c.Length - 16).Text = "to audio"
Set-ShapeGeometryEmu $audioBox $null $null $null 1169551

# ---------------------------------------------------------------------------
# 8) Connector feeding the audio box: shift down and shrink to match
# ---------------------------------------------------------------------------

$audioConn = Get-ShapeByName $s "Straight Arrow Connector 98"
Set-ShapeGeometryEmu $audioConn $null 4135682 $null 696808
